$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2375690607734807
$ws.Range("C2").Value = 0.4861878453038674
$ws.Range("J2").Value = 0.02209944751381215
$ws.Range("P2").Value = 0.138121546961326
$ws.Range("S2").Value = 0.1160220994475138
$ws.Range("B3").Value = 0.00546448087431694
$ws.Range("C3").Value = 0.04371584699453552
$ws.Range("J3").Value = 0.02185792349726776
$ws.Range("P3").Value = 0.7431693989071039
$ws.Range("S3").Value = 0.185792349726776
$ws.Range("J4").Value = 0.046875
$ws.Range("P4").Value = 0.59375
$ws.Range("S4").Value = 0.359375
$ws.Range("B6").Value = 0.09722222222222222
$ws.Range("D6").Value = 0.02777777777777778
$ws.Range("E6").Value = 0.004629629629629629
$ws.Range("F6").Value = 0.07407407407407407
$ws.Range("J6").Value = 0.2731481481481481
$ws.Range("O6").Value = 0.03240740740740741
$ws.Range("Q6").Value = 0.09722222222222222
$ws.Range("R6").Value = 0.07407407407407407
$ws.Range("S6").Value = 0.3194444444444444
$ws.Range("B7").Value = 0.0963302752293578
$ws.Range("D7").Value = 0.01834862385321101
$ws.Range("E7").Value = 0.004587155963302753
$ws.Range("F7").Value = 0.02752293577981652
$ws.Range("J7").Value = 0.1743119266055046
$ws.Range("O7").Value = 0.009174311926605505
$ws.Range("Q7").Value = 0.1605504587155963
$ws.Range("R7").Value = 0.08256880733944955
$ws.Range("S7").Value = 0.426605504587156
$ws.Range("B8").Value = 0.08490566037735849
$ws.Range("D8").Value = 0.0259433962264151
$ws.Range("F8").Value = 0.04716981132075472
$ws.Range("J8").Value = 0.1580188679245283
$ws.Range("O8").Value = 0.01179245283018868
$ws.Range("Q8").Value = 0.1863207547169811
$ws.Range("R8").Value = 0.1061320754716981
$ws.Range("S8").Value = 0.3797169811320755
$ws.Range("B9").Value = 0.08900523560209424
$ws.Range("D9").Value = 0.02617801047120419
$ws.Range("E9").Value = 0.005235602094240838
$ws.Range("F9").Value = 0.06806282722513089
$ws.Range("J9").Value = 0.1413612565445026
$ws.Range("O9").Value = 0.02094240837696335
$ws.Range("Q9").Value = 0.1675392670157068
$ws.Range("R9").Value = 0.1099476439790576
$ws.Range("S9").Value = 0.3717277486910995
$ws.Range("B10").Value = 0.1254307374224673
$ws.Range("D10").Value = 0.02756719503790489
$ws.Range("E10").Value = 0.002067539627842867
$ws.Range("F10").Value = 0.06547208821502412
$ws.Range("J10").Value = 0.129565816678153
$ws.Range("O10").Value = 0.01309441764300482
$ws.Range("Q10").Value = 0.2143349414197105
$ws.Range("R10").Value = 0.07236388697450034
$ws.Range("S10").Value = 0.3501033769813922
$ws.Range("G11").Value = 0.1295681063122923
$ws.Range("J11").Value = 0.08305647840531562
$ws.Range("K11").Value = 0.1926910299003322
$ws.Range("L11").Value = 0.5813953488372093
$ws.Range("S11").Value = 0.0132890365448505
$ws.Range("G12").Value = 0.7955801104972375
$ws.Range("J12").Value = 0.1712707182320442
$ws.Range("K12").Value = 0.005524861878453038
$ws.Range("L12").Value = 0.02209944751381215
$ws.Range("S12").Value = 0.005524861878453038
$ws.Range("G13").Value = 0.6721311475409836
$ws.Range("J13").Value = 0.2950819672131147
$ws.Range("S13").Value = 0.03278688524590164
$ws.Range("F15").Value = 0.01762114537444934
$ws.Range("H15").Value = 0.1585903083700441
$ws.Range("I15").Value = 0.05726872246696035
$ws.Range("J15").Value = 0.4096916299559472
$ws.Range("K15").Value = 0.07929515418502203
$ws.Range("M15").Value = 0.013215859030837
$ws.Range("O15").Value = 0.03964757709251102
$ws.Range("S15").Value = 0.2246696035242291
$ws.Range("F16").Value = 0.01333333333333333
$ws.Range("H16").Value = 0.1555555555555556
$ws.Range("I16").Value = 0.07111111111111111
$ws.Range("J16").Value = 0.4622222222222222
$ws.Range("K16").Value = 0.1022222222222222
$ws.Range("M16").Value = 0.02666666666666667
$ws.Range("O16").Value = 0.07111111111111111
$ws.Range("S16").Value = 0.09777777777777778
$ws.Range("F17").Value = 0.01691331923890063
$ws.Range("H17").Value = 0.1627906976744186
$ws.Range("I17").Value = 0.08879492600422834
$ws.Range("J17").Value = 0.4418604651162791
$ws.Range("K17").Value = 0.07822410147991543
$ws.Range("M17").Value = 0.03382663847780127
$ws.Range("O17").Value = 0.0613107822410148
$ws.Range("S17").Value = 0.1162790697674419
$ws.Range("F18").Value = 0.01470588235294118
$ws.Range("H18").Value = 0.1372549019607843
$ws.Range("I18").Value = 0.07843137254901961
$ws.Range("J18").Value = 0.4803921568627451
$ws.Range("K18").Value = 0.107843137254902
$ws.Range("M18").Value = 0.0196078431372549
$ws.Range("O18").Value = 0.05882352941176471
$ws.Range("S18").Value = 0.1029411764705882
$ws.Range("F19").Value = 0.01549186676994578
$ws.Range("H19").Value = 0.1944229279628195
$ws.Range("I19").Value = 0.08365608055770721
$ws.Range("J19").Value = 0.3865220759101472
$ws.Range("K19").Value = 0.1037955073586367
$ws.Range("M19").Value = 0.02633617350890782
$ws.Range("O19").Value = 0.07281177381874517
$ws.Range("S19").Value = 0.1169635941130906

Write-Output "Applied 110 cell updates"